# Apply the 2024-03-23 cryptos-list refresh (prices / 1h-volume / two rank swaps).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.427.71"
$ws.Range("E2").Value = "  +2.66%  "
$ws.Range("D3").Value = "3.408.90"
$ws.Range("E3").Value = "  +2.02%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.630"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.91%  "
$ws.Range("D8").Value = "3.401.93"
$ws.Range("E8").Value = "  +2.09%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.171"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +12.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.634"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.01"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000280"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.23%  "
$ws.Range("D15").Value = "3.947.13"
$ws.Range("E15").Value = "  +1.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.71%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.119"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.10%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.401.01"
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("D19").Value = "65.430.83"
$ws.Range("E19").Value = "  +2.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "473.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +15.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +16.60%  "
$ws.Range("E24").Value = "  +3.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "87.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.50"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.89"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.56"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "62.60"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "575.29"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.45%  "
$ws.Range("E35").Value = "  +2.15%  "
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.141"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.36%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.53"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.55%  "
$ws.Range("D40").Value = "0.0₃0762"
$ws.Range("E40").Value = "  +3.24%  "
$ws.Range("E41").Value = "  +2.36%  "
$ws.Range("D42").Value = "3.096.35"
$ws.Range("E42").Value = "  -1.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("E44").Value = "  +2.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0419"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.85%  "
$ws.Range("E46").Value = "  +3.77%  "
$ws.Range("E47").Value = "  +5.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.87%  "
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "137.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.84%  "
$ws.Range("E51").Value = "  +3.42%  "

Write-Host "cryptos list updated"
